# add case 2 3,5dof actio & reset
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.141, -0.08699999999999999, 0.779, 96, 22, 47),
    @(-0.122, -0.127, 0.205, -90, -53, -122),
    @(0.801, 0.318, 0.398, 89, -23, -141),
    @(-0.205, -0.034, -0.755, -133, 1, 121),
    @(-0.239, -0.045, 0.769, -139, 19, 131),
    @(0.602, -0.353, -0.484, 162, -20, 94),
    @(0.405, 0.143, 0.166, -24, 10, -77),
    @(-0.119, 0.06, 0.278, -143, -14, 128),
    @(-0.209, 0.032, 1.071, 60, 63, -134),
    @(-0.376, 0.273, 0.922, 97, 37, -55)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
